$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

# Fix existing row 16: columns B..K currently blank -> literal "nan" text
# (matches the pattern already used by earlier rows 14 and 15)
foreach ($col in @("B","C","D","E","F","G","H","I","J","K")) {
    $ws.Range($col + "16").Value = "nan"
}

# Append a new service-event row (row 17) for Card12
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "12"
$ws.Range("A17").Style = "Normal"

$ws.Range("L17").Value = "19/8/2025"
$ws.Range("M17").Value = "فني"
$ws.Range("N17").Value = "قطع سير كويلر مسنن 1270"
$ws.Range("O17").Value = "تم تغير سير 1270"
